# "Generate Report for Archive"
#
# Change every "Status" value that currently reads "Ready for handoff" to
# "In Translation" (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3) and shrink the
# now-narrower "Status" columns to fit the shorter text (Overview columns
# E/F, and column C on the zh-cn / de-de sheets).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$newStatus = "In Translation"

# --- Overview sheet: zh-cn status in column E, de-de status in column F ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

# --- zh-cn / de-de sheets: Status is column C ---
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# --- Re-fit the Status columns now that the text is shorter ---
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5
$zhcn.Columns.Item(3).ColumnWidth = 12.5
$dede.Columns.Item(3).ColumnWidth = 12.5
